$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Log a new volunteer session: 9:30PM -> 10:45PM on 1-20-2018 (75 minutes),
# entered on the next free row (25) right above the existing, already
# formula-driven "Total Project Hours:" row (29, with a formula summing
# C2:C28 so the new row is automatically picked up).
$ws.Range("A25").Value = "9:30PM 1-20-2018"
$ws.Range("B25").Value = "10:45PM 1-20-2018"
$ws.Range("B25").NumberFormat = "HH:MM:SS\ AM/PM"
$ws.Range("C25").Value = 75

# Update the view state to match the recorded selection/scroll position.
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("B26").Select() | Out-Null
